$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (A:T). Row 1 header is unchanged.
$data = @(
    @("ECs","Sfrp1","Fzd6","ECs",3,1,0.3669776666666666,1.100933,0.0181959334720815,0.0181959334720815,2,0.6666666666666666,6.220104333333333,18.660313,0.9231675752574263,0.9231675752574263,2.282639374669889,20.543754372029,0.01679789578296692,0.01679789578296692),
    @("ECs","Sfrp1","Fzd6","FAPs",3,1,0.3669776666666666,1.100933,0.0181959334720815,0.0181959334720815,2,0.6666666666666666,0.09168666666666665,0.27506,0.01360783569119702,0.01360783569119702,0.03364695899777777,0.30282263098,0.0002476072729360371,0.0002476072729360371),
    @("ECs","Sfrp1","Fzd6","sCs",3,1,0.3669776666666666,1.100933,0.0181959334720815,0.0181959334720815,3,1,0.4259936666666667,1.277981,0.06322458905137664,0.06322458905137664,0.1563301618081111,1.406971456273,0.001150430416178542,0.001150430416178542),
    @("FAPs","Sfrp1","Fzd6","ECs",3,1,17.23456,51.70368000000001,0.8545449373774706,0.8545449373774706,2,0.6666666666666666,6.220104333333333,18.660313,0.9231675752574263,0.9231675752574263,107.2007613390933,964.80685205184,0.7888881777872687,0.7888881777872687),
    @("FAPs","Sfrp1","Fzd6","FAPs",3,1,17.23456,51.70368000000001,0.8545449373774706,0.8545449373774706,2,0.6666666666666666,0.09168666666666665,0.27506,0.01360783569119702,0.01360783569119702,1.580179357866667,14.2216142208,0.01162850709857686,0.01162850709857686),
    @("FAPs","Sfrp1","Fzd6","sCs",3,1,17.23456,51.70368000000001,0.8545449373774706,0.8545449373774706,3,1,0.4259936666666667,1.277981,0.06322458905137664,0.06322458905137664,7.341813407786668,66.07632067008001,0.05402825249162496,0.05402825249162496),
    @("sCs","Sfrp1","Fzd6","ECs",3,1,2.566576666666667,7.69973,0.1272591291504479,0.1272591291504479,2,0.6666666666666666,6.220104333333333,18.660313,0.9231675752574263,0.9231675752574263,15.96437464616555,143.67937181549,0.1174815016871907,0.1174815016871907),
    @("sCs","Sfrp1","Fzd6","FAPs",3,1,2.566576666666667,7.69973,0.1272591291504479,0.1272591291504479,2,0.6666666666666666,0.09168666666666665,0.27506,0.01360783569119702,0.01360783569119702,0.2353208593111111,2.1178877338,0.001731721319684116,0.001731721319684116),
    @("sCs","Sfrp1","Fzd6","sCs",3,1,2.566576666666667,7.69973,0.1272591291504479,0.1272591291504479,3,1,0.4259936666666667,1.277981,0.06322458905137664,0.06322458905137664,1.093345405014444,9.84010864513,0.008045906143573136,0.008045906143573135)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}
